# Remove all the paragraphs after "Please describe the advantages of your
# proposed technology solution..." (the 5th paragraph) through to the end
# of the body content, leaving the final section properties (sectPr)
# intact.
$d = $word.ActiveDocument

$startPara = $d.Paragraphs.Item(6)
$startRange = $startPara.Range.Start
$endRange = $d.Content.End

$r = $d.Range($startRange, $endRange)
$r.Delete()
